$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Mark Yean Tuck Ming'
$ws.Range("B2").Value = '+6012-225 1051'
$ws.Range("C2").Value = 'yeanmark@gmail.com'
$ws.Range("F2").Value = '[{''job_title'': ''Data Engineer'', ''job_company'': ''CelcomDigi Sdn Bhd'', ''Industries'': [''Telecommunications''], ''start_date'': ''2022-05'', ''end_date'': ''2023-05'', ''job_location'': ''N/A'', ''job_duration'': ''1.0''}, {''job_title'': ''Data Analyst Intern'', ''job_company'': ''MC Crenergy Sdn Bhd'', ''Industries'': [''Unknown''], ''start_date'': ''2019-11'', ''end_date'': ''2020-03'', ''job_location'': ''N/A'', ''job_duration'': ''0.3''}]'
$ws.Range("G2").Value = '[{''Country'': ''Malaysia'', ''State'': ''Selangor'', ''City'': ''Puchong''}]'
$ws.Range("H2").Value = '[{''field_of_study'': ''Intelligence System'', ''level'': ''Bachelor Degree'', ''cgpa'': ''N/A'', ''university'': ''Asia Pacific University'', ''start_date'': ''2018-02'', ''year_of_graduation'': ''2021''}, {''field_of_study'': ''Information Technology'', ''level'': ''Foundation'', ''cgpa'': ''N/A'', ''university'': ''Asia Pacific University'', ''start_date'': ''2017-02'', ''year_of_graduation'': ''2017''}]'
$ws.Range("I2").Value = '[''Project Management Workshop Series'', ''Body Language Secrets: Nonverbal Communication Strategies'', ''Computer Hacking Forensic Investigator'', ''Certified Ethical Hacking Ver. 10'', ''APIIT Certified Security Professional'', ''Premier-Pride Challenge 2019'', ''Artificial Intelligence Modules'', ''AI and Machine Learning Competence for Industry 4.0'', ''Certified Engineer in Computer Vision'', ''Dataiku Core Designer Certificate'']'
$ws.Range("J2").Value = '[''Python'', ''Selenium'', ''BeautifulSoup'', ''Pandas'', ''Numpy'', ''Matplotlib'', ''Seaborn'', ''Plotly'', ''Supervised ML models'', ''Unsupervised ML models'', ''Natural Language Process'', ''Computer Vision'', ''Flask'', ''DJango'', ''Dash'', ''MySQL'', ''Teradata SQL'', ''Oracle Data Integrator'', ''Linux'', ''C++'', ''C#'', ''Java'', ''Microsoft Power BI'', ''OLAP Cube'', ''SSIS'', ''ASP .net'', ''DevExpress'']'
$ws.Range("K2").Value = '[''Chinese(Mandarin)'', ''Chinese(Cantonese)'', ''English'', ''Malay'']'
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 3.2
$ws.Range("N2").Value = 6.83
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3
$ws.Range("T2").Value = 8
$ws.Range("Z2").Value = '**Alignments:**
1. **Education Background:** The candidate has a Bachelor''s degree in Intelligence System, which aligns with the job description that may require a relevant educational background.
2. **Skills:** The candidate possesses a wide range of technical skills such as Python, Selenium, Pandas, Numpy, etc., which align with the technical skills mentioned in the job requirements.
3. **Previous Job Roles:** The candidate has experience as a Data Engineer and Data Analyst Intern, which aligns with the job description that may require experience in similar roles.
4. **Professional Certificates:** The candidate holds certificates related to Project Management, Artificial Intelligence, and Machine Learning, which align with the job requirements that may require relevant certifications.
5. **Languages:** The candidate is proficient in English, which aligns with the job requirements that may require proficiency in English.
**Misalignments:**
1. **Education Background:** The candidate''s foundation in Information Technology may not directly align with the job description''s requirement for a specific field of study.
2. **Skills:** While the candidate has a diverse set of technical skills, some of the skills mentioned may not be directly relevant to the job requirements, such as C++, C#, and ASP .net.
3. **Previous Job Roles:** The candidate''s previous job roles are in the Telecommunications and Unknown industries, which may not directly align with the industry or domain specified in the job description.
4. **Professional Certificates:** Some of the certificates held by the candidate, such as Body Language Secrets and Certified Ethical Hacking, may not directly align with the job requirements unless explicitly mentioned in the job description.
5. **Languages:** The candidate''s proficiency in Chinese (Mandarin and Cantonese) and Malay may not be directly relevant to the job requirements unless specified in the job description.'
$ws.Range("AA2").Value = 28.03

$ws.Rows(2).AutoFit()
